$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for the season record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (style) from an existing header cell (A1) onto
# the new header cells so they match the look of the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the season record (Wins=78, Losses=84, Ties=0) for every player row
for ($r = 2; $r -le 58; $r++) {
    $ws.Cells.Item($r, 30).Value = 78
    $ws.Cells.Item($r, 31).Value = 84
    $ws.Cells.Item($r, 32).Value = 0
}
